$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.072.54"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "2.301.46"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("D5").Value = "'310.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'100.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.11%  "
$ws.Range("D7").Value = "'0.537"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.79%  "
$ws.Range("D10").Value = "'36.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  +7.33%  "
$ws.Range("D14").Value = "2.657.72"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'14.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.52%  "
$ws.Range("D16").Value = "2.300.33"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("D18").Value = "42.995.62"
$ws.Range("E18").Value = "  +2.60%  "
$ws.Range("D19").Value = "'12.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "'6.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Value = "'68.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'240.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").Value = "'2.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'24.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "'38.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.25%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "'167.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +3.51%  "
$ws.Range("D36").Value = "'0.0739"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("E40").Value = "  +0.74%  "
$ws.Range("D41").Value = "'4.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.86%  "
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0289"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "1.971.18"
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("D46").Value = "'3.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("D47").Value = "'9.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "'55.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").Value = "'2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +16.52%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.527.37"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
